$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investments")

# Update sector max-weight constraints (C5:C15) from 15% to 10%
$ws.Range("C5:C15").Value = 0.1

# Update the selection to match the saved view state (C5:C15, active cell C5)
$ws.Range("C5:C15").Select()
